$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("# of RMAs") before the existing Quantity column,
# shifting the old Quantity column to C.
$ws.Range("B1").EntireColumn.Insert()

$ws.Range("B1").Value = "# of RMAs"
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 65
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 4
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 3

$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()

$ws.Range("A2:C9").Select() | Out-Null
